# Add application details for Wharton (row 8) on the "Application Details" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Application Details")

# Deadline (Dec 15, 2023), shown as "d-mmm" like other deadline cells.
$ws.Range("C8").NumberFormat = "d-mmm"
$ws.Range("C8").Value = 45275

# Fee (same currency-ish number format as the other Fee cells in this column).
$ws.Range("E8").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E8").Value = 80

# Resume? / SoP Notes / Writing sample length -- all still TBD.
$ws.Range("F8").Value = "TBD"
$ws.Range("G8").Value = "TBD"
$ws.Range("H8").Value = "TBD"

# # letters of rec
$ws.Range("I8").Value = 3

# Other requirements?
$ws.Range("J8").Value = "Additional optional essays"
